$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '29.023.67'
$ws.Range('E2').Value = '  -0.74%  '
Set-TextValue $ws.Range('D3') '1.831.44'
$ws.Range('E3').Value = '  -0.72%  '
Set-TextValue $ws.Range('D4') '0.9993'
$ws.Range('E4').Value = '  -0.03%  '
Set-TextValue $ws.Range('D5') '242.34'
$ws.Range('E5').Value = '  -0.23%  '
Set-TextValue $ws.Range('D6') '0.6252'
$ws.Range('E6').Value = '  -5.84%  '
Set-TextValue $ws.Range('D7') '0.9999'
$ws.Range('E7').Value = '  -0.02%  '
Set-TextValue $ws.Range('D8') '0.07556'
$ws.Range('E8').Value = '  +1.16%  '
Set-TextValue $ws.Range('D9') '0.2919'
$ws.Range('E9').Value = '  -1.53%  '
Set-TextValue $ws.Range('D10') '22.54'
$ws.Range('E10').Value = '  -3.40%  '
Set-TextValue $ws.Range('D11') '0.07728'
$ws.Range('E11').Value = '  -0.40%  '
Set-TextValue $ws.Range('D12') '1.833.36'
$ws.Range('E12').Value = '  -0.74%  '
Set-TextValue $ws.Range('D13') '4.946'
$ws.Range('E13').Value = '  -1.62%  '
Set-TextValue $ws.Range('D14') '0.6645'
$ws.Range('E14').Value = '  -1.52%  '
$ws.Range('E15').Value = '  +14.66%  '
Set-TextValue $ws.Range('D16') '82.70'
$ws.Range('E16').Value = '  -0.88%  '
Set-TextValue $ws.Range('D17') '6.030'
$ws.Range('E17').Value = '  -2.59%  '
Set-TextValue $ws.Range('D18') '28.987.77'
$ws.Range('E18').Value = '  -0.73%  '
Set-TextValue $ws.Range('D19') '226.84'
$ws.Range('E19').Value = '  -0.28%  '
Set-TextValue $ws.Range('D20') '12.32'
$ws.Range('E20').Value = '  -1.76%  '
Set-TextValue $ws.Range('D21') '0.9992'
$ws.Range('E21').Value = '  -0.11%  '
Set-TextValue $ws.Range('D22') '7.160'
$ws.Range('E22').Value = '  -0.90%  '
Set-TextValue $ws.Range('D23') '1.0000'
$ws.Range('E23').Value = '  -0.04%  '
Set-TextValue $ws.Range('D24') '158.06'
$ws.Range('E24').Value = '  -0.49%  '
Set-TextValue $ws.Range('D25') '8.465'
$ws.Range('E25').Value = '  -2.03%  '
Set-TextValue $ws.Range('D26') '0.1373'
$ws.Range('E26').Value = '  -2.04%  '
Set-TextValue $ws.Range('D27') '17.93'
$ws.Range('E27').Value = '  -0.70%  '
Set-TextValue $ws.Range('D28') '1.488'
$ws.Range('E28').Value = '  -1.81%  '
Set-TextValue $ws.Range('D29') '4.089'
$ws.Range('E29').Value = '  -1.32%  '
Set-TextValue $ws.Range('D30') '4.018'
$ws.Range('E30').Value = '  -0.82%  '
Set-TextValue $ws.Range('D31') '1.195'
$ws.Range('E31').Value = '  -0.52%  '
Set-TextValue $ws.Range('D32') '0.05198'
$ws.Range('E32').Value = '  -3.38%  '
Set-TextValue $ws.Range('D33') '1.847'
$ws.Range('E33').Value = '  -0.16%  '
Set-TextValue $ws.Range('D34') '0.7363'
$ws.Range('E34').Value = '  -1.51%  '
Set-TextValue $ws.Range('D35') '1.139'
$ws.Range('E35').Value = '  -2.04%  '
Set-TextValue $ws.Range('D36') '2.697'
$ws.Range('E36').Value = '  +1.88%  '
Set-TextValue $ws.Range('D37') '1.245.91'
$ws.Range('E37').Value = '  -4.09%  '
Set-TextValue $ws.Range('D38') '2.761'
$ws.Range('E38').Value = '  +0.31%  '
Set-TextValue $ws.Range('D39') '0.01784'
$ws.Range('E39').Value = '  -0.83%  '
Set-TextValue $ws.Range('D40') '6.332'
$ws.Range('E40').Value = '  -0.85%  '
Set-TextValue $ws.Range('D41') '0.8953'
$ws.Range('E41').Value = '  -1.26%  '
Set-TextValue $ws.Range('D42') '1.000'
$ws.Range('E42').Value = '  +0.02%  '
Set-TextValue $ws.Range('D43') '101.32'
$ws.Range('E43').Value = '  -2.25%  '
Set-TextValue $ws.Range('D44') '1.979.24'
$ws.Range('E44').Value = '  -0.73%  '
Set-TextValue $ws.Range('D45') '0.00000000123'
$ws.Range('E45').Value = '  +1.57%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D46') '64.02'
$ws.Range('E46').Value = '  -2.07%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D47') '0.5111'
$ws.Range('E47').Value = '  -0.59%  '
Set-TextValue $ws.Range('D48') '0.4039'
$ws.Range('E48').Value = '  +0.28%  '
Set-TextValue $ws.Range('D49') '8.853'
$ws.Range('E49').Value = '  +0.37%  '
Set-TextValue $ws.Range('D50') '0.05753'
$ws.Range('E50').Value = '  -1.89%  '
Set-TextValue $ws.Range('D51') '1.637'
$ws.Range('E51').Value = '  -6.78%  '
